$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4114.784832046243,
    4114.784832046243,
    4114.784832046243,
    4114.784832046243,
    4114.784832046243,
    4044.468177775046,
    4032.430749483882,
    4032.430749483882,
    3975.49186220928,
    3946.290375547124,
    3906.852602915739
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
